$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.76995849609375
$ws.Range("B1").Value = 5.384731769561768
$ws.Range("C1").Value = 6.220949172973633
$ws.Range("D1").Value = 6.496116638183594
$ws.Range("E1").Value = 4.951488494873047
